$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the total marks error on the "Marking" row (row 11)
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Correct the total marks error on the "Total" row (row 12)
$ws.Range("B12").Value = 92
$ws.Range("C12").Value = -10
$ws.Range("E12").Value = "82 / 112"
